# Generate Report for Handback
#
# For the "a9406554-f0a5-4a8f-9962-da3481f9a657" source file, a handback
# file showed up whose content/version didn't match the latest source
# doc. Record that on both the zh-cn and de-de report sheets:
#   - fill in the (previously empty) "Latest Target File" / "Latest
#     Handback File" / "Latest Handback DateTime" cells for row 7
#   - note the version mismatch in "Error Detail"
#   - widen the "Error Detail" column so the message is readable

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d3ac186fb989ea7a80b81324c1d56f1b66abe9de/e2e/a9406554-f0a5-4a8f-9962-da3481f9a657.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/218f5dcce253236be934d8cecded295433d17323/e2e/a9406554-f0a5-4a8f-9962-da3481f9a657.md."
$handbackMdTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/218f5dcce253236be934d8cecded295433d17323/e2e/a9406554-f0a5-4a8f-9962-da3481f9a657.md"
$handbackMdDisplay = "a9406554-f0a5-4a8f-9962-da3481f9a657.md"

function Update-HandbackRow7 {
    param(
        [string]$SheetName,
        [string]$HandbackDateTime,
        [string]$HandbackXlf
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # "Latest Target File" (I7): link to the handback markdown, same way
    # the other rows in column A already do.
    $ws.Hyperlinks.Add($ws.Range("I7"), $handbackMdTarget, "", "", $handbackMdDisplay)

    # "Latest Handback File" (J7)
    $ws.Range("J7").Value = $HandbackXlf

    # "Latest Handback DateTime" (K7)
    $ws.Range("K7").Value = $HandbackDateTime

    # "Error Detail" (P7)
    $ws.Range("P7").Value = $errorDetail

    # Widen the Error Detail column (16 / P) so the message is legible.
    $ws.Columns.Item(16).ColumnWidth = 39.17
}

Update-HandbackRow7 -SheetName "zh-cn" -HandbackDateTime "2016-08-28 14:42:43" -HandbackXlf "a9406554-f0a5-4a8f-9962-da3481f9a657.08e1a2f228282f04b2843210c52db256bac06941.zh-cn.xlf"
Update-HandbackRow7 -SheetName "de-de" -HandbackDateTime "2016-08-28 14:42:50" -HandbackXlf "a9406554-f0a5-4a8f-9962-da3481f9a657.08e1a2f228282f04b2843210c52db256bac06941.de-de.xlf"
